$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value to J8, matching the style already used by the other
# numeric cells in that row/column (same xf as J7/I8 etc.)
$ws.Range("J8").Value = 1450
$ws.Range("J7").Copy()
$ws.Range("J8").PasteSpecial(-4122)  # xlPasteFormats

# Update the selection to L8 (as seen in the diff)
$ws.Range("L8").Select()
